# Scheduled runner update: refresh computed pricing/profit figures across
# several item rows on each crafting-class sheet (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 532.9
$ws.Range("J17").Value = 532.9
$ws.Range("L17").Value = 1598.7
$ws.Range("N17").Value = -1934.7

$ws.Range("H29").Value = 1470.7097
$ws.Range("I29").Value = 1699.75
$ws.Range("J29").Value = 1436.7778
$ws.Range("K29").Value = 5099.25
$ws.Range("L29").Value = 4310.3334
$ws.Range("M29").Value = -4818.25
$ws.Range("N29").Value = -4872.3334

$ws.Range("H58").Value = 291.66666
$ws.Range("I58").Value = 78.125
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 234.375
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -84.375
$ws.Range("N58").Value = -6300

$ws.Range("H132").Value = 8720.853999999999
$ws.Range("I132").Value = 2114.606
$ws.Range("K132").Value = 6343.818000000001
$ws.Range("M132").Value = -3813.818000000001

$ws.Range("H137").Value = 7755497.5
$ws.Range("I137").Value = 1159.6364
$ws.Range("J137").Value = 15879090
$ws.Range("K137").Value = 3478.9092
$ws.Range("L137").Value = 47637270
$ws.Range("M137").Value = -928.9092000000001
$ws.Range("N137").Value = -47642370

$ws.Range("H138").Value = 5627.8394
$ws.Range("J138").Value = 7880.5264
$ws.Range("L138").Value = 23641.5792
$ws.Range("N138").Value = -33921.5792

$ws.Range("H140").Value = 58421.547
$ws.Range("J140").Value = 57192.8
$ws.Range("L140").Value = 57192.8
$ws.Range("N140").Value = -67552.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6150.478
$ws.Range("I32").Value = 3530.4849
$ws.Range("J32").Value = 12801.23
$ws.Range("K32").Value = 3530.4849
$ws.Range("L32").Value = 12801.23
$ws.Range("M32").Value = -3243.4849
$ws.Range("N32").Value = -13375.23

$ws.Range("H63").Value = 4075.25
$ws.Range("J63").Value = 3806
$ws.Range("L63").Value = 3806
$ws.Range("N63").Value = -5178

$ws.Range("H66").Value = 4075.25
$ws.Range("J66").Value = 3806
$ws.Range("L66").Value = 19030
$ws.Range("N66").Value = -25894

$ws.Range("H74").Value = 7577285
$ws.Range("I74").Value = 11905720
$ws.Range("J74").Value = 2524.5
$ws.Range("K74").Value = 11905720
$ws.Range("L74").Value = 2524.5
$ws.Range("M74").Value = -11904846
$ws.Range("N74").Value = -4272.5

$ws.Range("H77").Value = 7577285
$ws.Range("I77").Value = 11905720
$ws.Range("J77").Value = 2524.5
$ws.Range("K77").Value = 59528600
$ws.Range("L77").Value = 12622.5
$ws.Range("M77").Value = -59524232
$ws.Range("N77").Value = -21358.5

$ws.Range("H92").Value = 15020275
$ws.Range("J92").Value = 15020275
$ws.Range("L92").Value = 15020275
$ws.Range("N92").Value = -15025267

$ws.Range("H97").Value = 601.5333000000001
$ws.Range("I97").Value = 608.7857
$ws.Range("K97").Value = 608.7857
$ws.Range("M97").Value = -112.7857

$ws.Range("H132").Value = 19437.486
$ws.Range("I132").Value = 20961.207
$ws.Range("K132").Value = 62883.621
$ws.Range("M132").Value = -60353.621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3582.5833
$ws.Range("I20").Value = 3510.3333
$ws.Range("J20").Value = 3799.3333
$ws.Range("K20").Value = 3510.3333
$ws.Range("L20").Value = 3799.3333
$ws.Range("M20").Value = -3263.3333
$ws.Range("N20").Value = -4293.3333

$ws.Range("H94").Value = 596262.5600000001
$ws.Range("I94").Value = 623274.5
$ws.Range("K94").Value = 623274.5
$ws.Range("M94").Value = -622823.5

$ws.Range("H134").Value = 8214
$ws.Range("I134").Value = 5537.25
$ws.Range("J134").Value = 9998.5
$ws.Range("K134").Value = 16611.75
$ws.Range("L134").Value = 29995.5
$ws.Range("M134").Value = -14076.75
$ws.Range("N134").Value = -35065.5

$ws.Range("H138").Value = 82500
$ws.Range("J138").Value = 82500
$ws.Range("L138").Value = 82500
$ws.Range("N138").Value = -92780

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 125000
$ws.Range("J141").Value = 125000
$ws.Range("L141").Value = 125000
$ws.Range("N141").Value = -135360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2255.5
$ws.Range("I16").Value = 2007.3334
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 2007.3334
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1720.3334
$ws.Range("N16").Value = -3574

$ws.Range("H31").Value = 1561.5
$ws.Range("I31").Value = 1006.6875
$ws.Range("K31").Value = 1006.6875
$ws.Range("M31").Value = -711.6875

$ws.Range("H34").Value = 1561.5
$ws.Range("I34").Value = 1006.6875
$ws.Range("K34").Value = 1006.6875
$ws.Range("M34").Value = -804.6875

$ws.Range("H94").Value = 2749.087
$ws.Range("J94").Value = 2695.5715
$ws.Range("L94").Value = 2695.5715
$ws.Range("N94").Value = -3597.5715

$ws.Range("H113").Value = 2255.5
$ws.Range("I113").Value = 2007.3334
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2007.3334
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 162.6666
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 978.17645
$ws.Range("I5").Value = 890.5
$ws.Range("K5").Value = 2671.5
$ws.Range("M5").Value = -2559.5

$ws.Range("H131").Value = 2857.319
$ws.Range("I131").Value = 914
$ws.Range("J131").Value = 2943.689
$ws.Range("K131").Value = 2742
$ws.Range("L131").Value = 8831.066999999999
$ws.Range("M131").Value = 2298
$ws.Range("N131").Value = -18911.067

$ws.Range("H135").Value = 978.17645
$ws.Range("I135").Value = 890.5
$ws.Range("K135").Value = 8014.5
$ws.Range("M135").Value = -5479.5

$ws.Range("H137").Value = 8474080
$ws.Range("J137").Value = 15736571
$ws.Range("L137").Value = 47209713
$ws.Range("N137").Value = -47219913

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 584683.9
$ws.Range("I80").Value = 1389546.4
$ws.Range("J80").Value = 16545.588
$ws.Range("K80").Value = 1389546.4
$ws.Range("L80").Value = 16545.588
$ws.Range("M80").Value = -1388548.4
$ws.Range("N80").Value = -18541.588

$ws.Range("H83").Value = 584683.9
$ws.Range("I83").Value = 1389546.4
$ws.Range("J83").Value = 16545.588
$ws.Range("K83").Value = 6947732
$ws.Range("L83").Value = 82727.94
$ws.Range("M83").Value = -6942740
$ws.Range("N83").Value = -92711.94

$ws.Range("H126").Value = 3767.8928
$ws.Range("I126").Value = 3069.9375
$ws.Range("K126").Value = 9209.8125
$ws.Range("M126").Value = -6739.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2323.9092
$ws.Range("I100").Value = 1840.3334
$ws.Range("K100").Value = 1840.3334
$ws.Range("M100").Value = -1299.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 21811.223
$ws.Range("I62").Value = 13100
$ws.Range("J62").Value = 26166.834
$ws.Range("K62").Value = 13100
$ws.Range("L62").Value = 26166.834
$ws.Range("M62").Value = -12476
$ws.Range("N62").Value = -27414.834

$ws.Range("H65").Value = 21811.223
$ws.Range("I65").Value = 13100
$ws.Range("J65").Value = 26166.834
$ws.Range("K65").Value = 65500
$ws.Range("L65").Value = 130834.17
$ws.Range("M65").Value = -62380
$ws.Range("N65").Value = -137074.17

$ws.Range("H132").Value = 30866426
$ws.Range("I132").Value = 4631257
$ws.Range("J132").Value = 83336770
$ws.Range("K132").Value = 13893771
$ws.Range("L132").Value = 250010310
$ws.Range("M132").Value = -13891241
$ws.Range("N132").Value = -250015370

$ws.Range("H136").Value = 8812.513999999999
$ws.Range("I136").Value = 3627
$ws.Range("K136").Value = 10881
$ws.Range("M136").Value = -8331
